$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new task row (row 11): task name in A11, status "Checked" in C11
$ws.Range("A11").Value = "Presentation"
$ws.Range("C11").Value = "Checked"

# Move active selection to D11 as in the final file
$ws.Range("D11").Select()
